$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StatusTableOverview")

# Update the "Code Edits" cell for the Philosophers row (row 2) with the new
# comment about the code edit, and give it the same "Good" style as the
# neighboring "Check in Afra" cell (D2).
$ws.Range("E2").Value = "msgsrv -> constructor, removed comments"
$ws.Range("E2").Style = "Good"

# Move the active selection to F15 (as recorded when the workbook was saved).
$ws.Range("F15").Select()
